# Update the "想去人数" (interested count) figures in column F on the
# "展览" and "全部类型" sheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$sheetExpo = $wb.Worksheets.Item("展览")
$sheetExpo.Range("F5").Value = 358
$sheetExpo.Range("F6").Value = 618
$sheetExpo.Range("F7").Value = 100
$sheetExpo.Range("F8").Value = 2072
$sheetExpo.Range("F9").Value = 10665
$sheetExpo.Range("F13").Value = 203
$sheetExpo.Range("F15").Value = 7525
$sheetExpo.Range("F20").Value = 3334

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F5").Value = 358
$sheetAll.Range("F6").Value = 618
$sheetAll.Range("F8").Value = 100
$sheetAll.Range("F9").Value = 2072
$sheetAll.Range("F12").Value = 10665
$sheetAll.Range("F16").Value = 203
$sheetAll.Range("F18").Value = 7525
$sheetAll.Range("F23").Value = 3334
